# Scheduled runner update: refresh market-price derived columns (H-N) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 382.375
$ws.Range("J17").Value = 382.375
$ws.Range("L17").Value = 1147.125
$ws.Range("N17").Value = -1483.125
$ws.Range("H40").Value = 2011.3334
$ws.Range("I40").Value = 2236.9
$ws.Range("J40").Value = 1806.2727
$ws.Range("K40").Value = 2236.9
$ws.Range("L40").Value = 1806.2727
$ws.Range("M40").Value = -2061.9
$ws.Range("N40").Value = -2156.2727
$ws.Range("H70").Value = 1538.2667
$ws.Range("I70").Value = 1409.1
$ws.Range("K70").Value = 4227.299999999999
$ws.Range("M70").Value = -3957.299999999999
$ws.Range("H73").Value = 1538.2667
$ws.Range("I73").Value = 1409.1
$ws.Range("K73").Value = 4227.299999999999
$ws.Range("M73").Value = -3291.299999999999
$ws.Range("H116").Value = 3205.4546
$ws.Range("I116").Value = 2903.524
$ws.Range("J116").Value = 3733.8333
$ws.Range("K116").Value = 2903.524
$ws.Range("L116").Value = 3733.8333
$ws.Range("M116").Value = 538.4760000000001
$ws.Range("N116").Value = -10617.8333
$ws.Range("H141").Value = 793.4286
$ws.Range("I141").Value = 793.4286
$ws.Range("K141").Value = 2380.2858
$ws.Range("M141").Value = 2799.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3826.6182
$ws.Range("I32").Value = 3728.9167
$ws.Range("K32").Value = 3728.9167
$ws.Range("M32").Value = -3441.9167
$ws.Range("H110").Value = 1187.35
$ws.Range("I110").Value = 938.1429000000001
$ws.Range("J110").Value = 1768.8334
$ws.Range("K110").Value = 938.1429000000001
$ws.Range("L110").Value = 1768.8334
$ws.Range("M110").Value = 1106.8571
$ws.Range("N110").Value = -5858.8334
$ws.Range("H132").Value = 2223.6875
$ws.Range("I132").Value = 1886.64
$ws.Range("K132").Value = 5659.92
$ws.Range("M132").Value = -3129.92

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2271.2856
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 2149.8333
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 2149.8333
$ws.Range("M20").Value = -2753
$ws.Range("N20").Value = -2643.8333
$ws.Range("H80").Value = 644.619
$ws.Range("I80").Value = 430
$ws.Range("K80").Value = 430
$ws.Range("M80").Value = 568
$ws.Range("H83").Value = 644.619
$ws.Range("I83").Value = 430
$ws.Range("K83").Value = 2150
$ws.Range("M83").Value = 2842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 47619884
$ws.Range("I16").Value = 50000788
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 50000788
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -50000501
$ws.Range("N16").Value = -2374
$ws.Range("H113").Value = 47619884
$ws.Range("I113").Value = 50000788
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 50000788
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -49998618
$ws.Range("N113").Value = -6140
$ws.Range("H134").Value = 1114.0322
$ws.Range("I134").Value = 929.0741
$ws.Range("J134").Value = 2362.5
$ws.Range("K134").Value = 2787.2223
$ws.Range("L134").Value = 7087.5
$ws.Range("M134").Value = -252.2223000000004
$ws.Range("N134").Value = -12157.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 101.08
$ws.Range("I12").Value = 156.42857
$ws.Range("J12").Value = 79.55556
$ws.Range("K12").Value = 469.28571
$ws.Range("L12").Value = 238.66668
$ws.Range("M12").Value = -296.28571
$ws.Range("N12").Value = -584.66668
$ws.Range("H32").Value = 2130.4546
$ws.Range("J32").Value = 2273.3
$ws.Range("L32").Value = 6819.900000000001
$ws.Range("N32").Value = -7385.900000000001
$ws.Range("H58").Value = 3067.5
$ws.Range("J58").Value = 3280
$ws.Range("L58").Value = 9840
$ws.Range("N58").Value = -10096
$ws.Range("H103").Value = 1474.4615
$ws.Range("I103").Value = 250
$ws.Range("J103").Value = 1841.8
$ws.Range("K103").Value = 750
$ws.Range("L103").Value = 5525.4
$ws.Range("M103").Value = 129
$ws.Range("N103").Value = -7283.4
$ws.Range("H131").Value = 14086877
$ws.Range("J131").Value = 2563.8
$ws.Range("L131").Value = 7691.400000000001
$ws.Range("N131").Value = -17771.4
$ws.Range("H134").Value = 3515.25
$ws.Range("I134").Value = 950.6
$ws.Range("J134").Value = 6079.9
$ws.Range("K134").Value = 2851.8
$ws.Range("L134").Value = 18239.7
$ws.Range("M134").Value = 2218.2
$ws.Range("N134").Value = -28379.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 112502250
$ws.Range("I70").Value = 83336330
$ws.Range("J70").Value = 200000000
$ws.Range("K70").Value = 83336330
$ws.Range("L70").Value = 200000000
$ws.Range("M70").Value = -83336060
$ws.Range("N70").Value = -200000540
$ws.Range("H73").Value = 112502250
$ws.Range("I73").Value = 83336330
$ws.Range("J73").Value = 200000000
$ws.Range("K73").Value = 83336330
$ws.Range("L73").Value = 200000000
$ws.Range("M73").Value = -83335394
$ws.Range("N73").Value = -200001872
$ws.Range("H80").Value = 6735.2
$ws.Range("I80").Value = 6707.4287
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 6707.4287
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -5709.4287
$ws.Range("N80").Value = -8796
$ws.Range("H83").Value = 6735.2
$ws.Range("I83").Value = 6707.4287
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 33537.14350000001
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -28545.14350000001
$ws.Range("N83").Value = -43984
$ws.Range("H102").Value = 2848.7097
$ws.Range("I102").Value = 3184.3076
$ws.Range("J102").Value = 2606.3333
$ws.Range("K102").Value = 3184.3076
$ws.Range("L102").Value = 2606.3333
$ws.Range("M102").Value = -1562.3076
$ws.Range("N102").Value = -5850.3333
$ws.Range("I107").Value = 837.13336
$ws.Range("J107").Value = 650.8333
$ws.Range("K107").Value = 837.13336
$ws.Range("L107").Value = 650.8333
$ws.Range("M107").Value = 1082.86664
$ws.Range("N107").Value = -4490.8333
$ws.Range("H132").Value = 1872.6428
$ws.Range("I132").Value = 1549.2632
$ws.Range("K132").Value = 4647.7896
$ws.Range("M132").Value = -2117.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1746.6666
$ws.Range("I82").Value = 1760
$ws.Range("J82").Value = 1733.3334
$ws.Range("K82").Value = 1760
$ws.Range("L82").Value = 1733.3334
$ws.Range("M82").Value = -1399
$ws.Range("N82").Value = -2455.3334
$ws.Range("H85").Value = 1746.6666
$ws.Range("I85").Value = 1760
$ws.Range("J85").Value = 1733.3334
$ws.Range("K85").Value = 1760
$ws.Range("L85").Value = 1733.3334
$ws.Range("M85").Value = -512
$ws.Range("N85").Value = -4229.3334
$ws.Range("H94").Value = 19249.5
$ws.Range("J94").Value = 19249.5
$ws.Range("L94").Value = 19249.5
$ws.Range("N94").Value = -20601.5
$ws.Range("H132").Value = 19974.723
$ws.Range("I132").Value = 964.4545000000001
$ws.Range("J132").Value = 49848
$ws.Range("K132").Value = 2893.3635
$ws.Range("L132").Value = 149544
$ws.Range("M132").Value = -363.3635000000004
$ws.Range("N132").Value = -154604

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 71434040
$ws.Range("I62").Value = 100006560
$ws.Range("K62").Value = 100006560
$ws.Range("M62").Value = -100005936
$ws.Range("H65").Value = 71434040
$ws.Range("I65").Value = 100006560
$ws.Range("K65").Value = 500032800
$ws.Range("M65").Value = -500029680
$ws.Range("H121").Value = 26400
$ws.Range("J121").Value = 26400
$ws.Range("L121").Value = 26400
$ws.Range("N121").Value = -29894
$ws.Range("H132").Value = 1463.9025
$ws.Range("J132").Value = 2198.9092
$ws.Range("L132").Value = 6596.7276
$ws.Range("N132").Value = -11656.7276
